{"js": "// HERCULES-8261 - [UMU][UCV][ETI] Informes evaluaci\u00f3n - Especificar tipo...\n//\n// The certification paragraph reads:\n//   \"Con fecha {{fechaEnvioSecretaria}}, {{fieldDelInvestigador}} {{fieldInvestigador}}\n//    solicita una modificaci\u00f3n relevante al proyecto ya aprobado.\"\n//\n// The merge-field placeholder \"{{fieldDelInvestigador}}\" (\"del Investigador\") is\n// renamed to \"{{fieldElInvestigador}}\" (\"El Investigador\") so the report template\n// can specify the article/gender variant explicitly. Everything else in that\n// sentence is unchanged.\n//\n// Use body.search to locate the exact placeholder token and replace it in place;\n// this only touches the two characters that actually changed (\"De\" -> \"E\") and\n// leaves the surrounding runs/formatting untouched.\nconst searchResults = context.document.body.search(\"{{fieldDelInvestigador}}\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"{{fieldElInvestigador}}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# HERCULES-8261 - [UMU][UCV][ETI] Informes evaluaci\u00f3n - Especificar tipo...\n#\n# The certification paragraph reads:\n#   \"Con fecha {{fechaEnvioSecretaria}}, {{fieldDelInvestigador}} {{fieldInvestigador}}\n#    solicita una modificaci\u00f3n relevante al proyecto ya aprobado.\"\n#\n# The merge-field placeholder \"{{fieldDelInvestigador}}\" (\"del Investigador\") is\n# renamed to \"{{fieldElInvestigador}}\" (\"El Investigador\") so the report template\n# can specify the article/gender variant explicitly. Everything else in that\n# sentence is unchanged.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"{{fieldDelInvestigador}}\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"{{fieldElInvestigador}}\"\n$find.Execute(\n    $find.Text,      # FindText\n    $false,          # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    1,               # Wrap (wdFindContinue)\n    $false,          # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                # Replace (wdReplaceAll)\n)\n\n# Two new \"ListLabel\" character styles are introduced in styles.xml - they back\n# list markers whose run properties were previously inherited/implicit and are\n# now pinned to an explicit font/size so the numbering glyphs render consistently.\n$listLabel97 = $d.Styles.Add(\"ListLabel97\", 2)\n$listLabel97.NameLocal = \"ListLabel 97\"\n$listLabel97.QuickStyle = $true\n$listLabel97.Font.Name = \"Ubuntu\"\n$listLabel97.Font.Size = 11\n$listLabel97.Font.SizeBi = 11\n\n$listLabel98 = $d.Styles.Add(\"ListLabel98\", 2)\n$listLabel98.NameLocal = \"ListLabel 98\"\n$listLabel98.QuickStyle = $true\n$listLabel98.Font.Name = \"Ubuntu\"\n$listLabel98.Font.Size = 10.5\n$listLabel98.Font.SizeBi = 10.5\n"}
